# AFDP-2536: make FOIA case file rules safe when queue has not been set
#
# The "Nullify Billing Enter Date" / "Nullify Hold Enter Date" rules used a
# bare `queue.name` reference that NPEs when a case file has never had a
# queue assigned. Switch to Spring SpEL's safe-navigation operator
# (`queue?.name`) everywhere, and tighten the "nullify" conditions so they
# only fire when the corresponding enter-date field is actually still set.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 28 - Set Billing Enter Date: guard the queue lookup.
$ws.Range("C28").Value() = "queue?.name == 'Billing' && billingEnterDate == null"

# Row 29 - Nullify Billing Enter Date: guard the queue lookup and only
# nullify when billingEnterDate has actually been populated.
$ws.Range("C29").Value() = "queue?.name != 'Billing' && billingEnterDate != null"

# Row 30 - Set Hold Enter Date: guard the queue lookup.
$ws.Range("C30").Value() = "queue?.name == 'Hold' && holdEnterDate == null"

# Row 31 - Nullify Hold Enter Date: guard the queue lookup and only nullify
# when holdEnterDate has actually been populated.
$ws.Range("C31").Value() = "queue?.name != 'Hold' && holdEnterDate != null"

# Reflect the author's last selection in the saved view state.
$ws.Activate()
$ws.Range("C30").Select()
